$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 238
$ws.Range("I2").Value = 603
$ws.Range("J2").Value = 2427
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 663
$ws.Range("M2").Value = 49
$ws.Range("N2").Value = 437
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 258
$ws.Range("T2").Value = 420
$ws.Range("V2").Value = 3732
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3824
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 46
$ws.Range("AA2").Value = 26
